$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<down>"
$ws.Range("C2").Value = 31

$ws.Range("B3").Value = "<callo>"
$ws.Range("C3").Value = 25

$ws.Range("B4").Value = "<down>"
$ws.Range("C4").Value = 32

$ws.Range("B5").Value = "<quebec>"
$ws.Range("C5").Value = 26

$ws.Range("C6").Value = 27

$ws.Range("B7").Value = "<otimd>"
$ws.Range("C7").Value = 30

$ws.Range("C8").Value = 26

$ws.Range("C9").Value = 33

$ws.Range("C10").Value = 30

$ws.Range("C11").Value = 24

$ws.Range("C12").Value = 21

$ws.Range("C13").Value = 30

$ws.Range("B14").Value = "<sie>"
$ws.Range("C14").Value = 30

$ws.Range("B15").Value = "<it>"
$ws.Range("C15").Value = 25

$ws.Range("C16").Value = 29

$ws.Range("B17").Value = "<sero>"
$ws.Range("C17").Value = 30

$ws.Range("C18").Value = 26
